# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.658.35"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.121.21"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.117.98"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.432"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "3.660.35"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "57.763.28"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "3.129.16"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.510"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0666"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "3.166.36"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.689"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "2.291.83"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.993"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.27%  "
